$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

# Give the new date cell (A7) the same formatting as the existing date
# column (A2:A6 -> bold/bordered/centered "date" style) before writing values.
$ws.Cells.Item(6, 1).Copy()
$ws.Cells.Item(7, 1).PasteSpecial(-4122)  # xlPasteFormats

# Append the new row of data after the existing last row (row 6 -> row 7)
$ws.Cells.Item(7, 1).Value = "06/05/2020 06:58:05"
$ws.Cells.Item(7, 2).Value = 6559.66
$ws.Cells.Item(7, 3).Value = 3447.35
$ws.Cells.Item(7, 4).Value = 10007.01
